$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.987.95"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.670.09"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.06"
$ws.Range("E5").Value = "  +4.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.94"
$ws.Range("E6").Value = "  +7.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.708"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.92"
$ws.Range("E10").Value = "  +10.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  -5.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000275"
$ws.Range("E12").Value = "  -5.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.22"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.263.97"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.674.08"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.96"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.12"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.898.70"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.51"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "401.25"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "87.85"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.96"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.58"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.68"
$ws.Range("E28").Value = "  -3.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.94"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.47"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "68.27"
$ws.Range("E32").Value = "  +5.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.37"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "45.27"
$ws.Range("E34").Value = "  +4.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "611.36"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.397"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  -14.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.135"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.92"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0427"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  -7.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.821.61"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.99"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.17"
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.55"
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  -2.50%  "
